# LOB1228.xlsx update — 2022 activation date + new English translations + revised "Programa" text
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Ativação:" value 01/01/2018 -> 01/01/2022 (kept as plain text, not a date) ---
# Going through a text formula + paste-as-values keeps the literal string
# "01/01/2022" instead of Excel auto-converting a bare date-looking value.
$ws.Range("B8").Formula = "=""01/01/2022"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)

$ws.Range("C8").Formula = "=""01/01/2022"""
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

# --- 2) New row 11 content under "Objectives:" (English objectives paragraph) ---
$ws.Range("B11").Value = "Provide students with knowledge about environmental management in companies, environmental policies, environmental management systems (EMS) and ISO 14000 series standards, enabling them to participate in the planning and implementation of an EMS in a company."
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C11").Value = "Provide students with knowledge about environmental management in companies, environmental policies, environmental management systems (EMS) and ISO 14000 series standards, enabling them to participate in the planning and implementation of an EMS in a company."
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- 3) New row 15 content under "Short syllabus:" ---
$ws.Range("B15").Value = "Environmental Management Systems; Iso 14000; Environmental Audit."
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").Value = "Environmental Management Systems; Iso 14000; Environmental Audit."
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# --- 4) "Programa:" body text (row 16) revised/expanded ---
$ws.Range("B16").Value = "Evolução das práticas de gestão ambiental empresarial;- Economia circular, conceitos e aplicações;- Responsabilidade social corporativa: conceito e programa;- Implantação do sistema de gerenciamento ambiental (SGA): conceitos e modelos;- Produção mais limpa;- Ferramentas de gestão focadas no produto;- Análise e otimização do ciclo de vida do produto;- Ecoinovação e Ecodesign;- Rotulagem ambiental;- Inovação e sustentabilidade;- Normas ISO 14001 (série ISO 14000), requisitos e orientações para uso e Certificações ambientais."
$ws.Range("C16").Value = "Evolução das práticas de gestão ambiental empresarial;- Economia circular, conceitos e aplicações;- Responsabilidade social corporativa: conceito e programa;- Implantação do sistema de gerenciamento ambiental (SGA): conceitos e modelos;- Produção mais limpa;- Ferramentas de gestão focadas no produto;- Análise e otimização do ciclo de vida do produto;- Ecoinovação e Ecodesign;- Rotulagem ambiental;- Inovação e sustentabilidade;- Normas ISO 14001 (série ISO 14000), requisitos e orientações para uso e Certificações ambientais."

# --- 5) New row 17 content under "Syllabus:" (English translation of Programa) ---
$ws.Range("B17").Value = "Evolution of corporate environmental management practices;- Circular economy, concepts and applications;- Corporate social responsibility: concept and program;- Implementation of the environmental management system (SGA): concepts and models;- Cleaner production;- Management tools focused on the product;- Analysis and optimization of the product life cycle;- Eco-innovation and Ecodesign;- Environmental labeling;- Innovation and sustainability;- ISO 14001 standards (ISO 14000 series), requirements and guidelines for use and Environmental Certifications."
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C17").Value = "Evolution of corporate environmental management practices;- Circular economy, concepts and applications;- Corporate social responsibility: concept and program;- Implementation of the environmental management system (SGA): concepts and models;- Cleaner production;- Management tools focused on the product;- Analysis and optimization of the product life cycle;- Eco-innovation and Ecodesign;- Environmental labeling;- Innovation and sustainability;- ISO 14001 standards (ISO 14000 series), requirements and guidelines for use and Environmental Certifications."
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
